$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column O (so old O..U shift to Q..W)
$ws.Range("O:P").Insert()

# Update header text for M1, N1 (same position, text changed)
$ws.Range("M1").Value = "Detected Predicates Doc Parent"
$ws.Range("N1").Value = "Detected Predicates Doc Related"

# Set headers for the two newly-inserted columns O1, P1
$ws.Range("O1").Value = "Correct Pred Predicates Parents"
$ws.Range("P1").Value = "Correct Pred Predicates Related"

# Populate the new O/P data columns (rows 2-6) with the recomputed values
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 2

$ws.Range("O3").Value = 3
$ws.Range("P3").Value = 3

$ws.Range("O4").Value = 2
$ws.Range("P4").Value = 2

$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0

$ws.Range("O6").Value = 2
$ws.Range("P6").Value = 2
